$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64
$ws.Range("A64").Value = 45932
$ws.Range("B64").Value = "四方坪站充电量(kw)"
$ws.Range("C64").Value = 797.86700000000008
$ws.Range("D64").Value = 1132.2600000000002
$ws.Range("E64").Value = 561.88099999999997
$ws.Range("F64").Value = 501.08199999999999
$ws.Range("G64").Value = 286.29700000000003
$ws.Range("H64").Value = 786.66499999999996
$ws.Range("I64").Value = 426.47900000000004
$ws.Range("J64").Value = 179.68899999999999
$ws.Range("K64").Value = 80.180000000000007
$ws.Range("L64").Value = 136.81299999999999
$ws.Range("M64").Value = 207.07200000000006
$ws.Range("N64").Value = 159.25
$ws.Range("O64").Value = 616.62499999999989
$ws.Range("P64").Value = 1192.2559999999999
$ws.Range("Q64").Value = 509.90899999999999
$ws.Range("R64").Value = 519.452
$ws.Range("S64").Value = 247.398
$ws.Range("T64").Value = 286.31099999999998
$ws.Range("U64").Value = 198.16
$ws.Range("V64").Value = 150.166
$ws.Range("W64").Value = 159.70700000000002
$ws.Range("X64").Value = 75.099999999999994
$ws.Range("Y64").Value = 130.12
$ws.Range("Z64").Value = 90.855000000000004

# Row 65
$ws.Range("A65").Value = 45932
$ws.Range("B65").Value = "高岭站充电量(kw)"
$ws.Range("C65").Value = 270.62299999999999
$ws.Range("D65").Value = 508.86800000000005
$ws.Range("E65").Value = 179.03199999999998
$ws.Range("F65").Value = 154.90299999999999
$ws.Range("G65").Value = 52.469000000000001
$ws.Range("H65").Value = 307.63000000000005
$ws.Range("I65").Value = 169.39099999999999
$ws.Range("J65").Value = 122.90700000000001
$ws.Range("K65").Value = 259.47200000000004
$ws.Range("L65").Value = 196.90699999999998
$ws.Range("M65").Value = 138.26
$ws.Range("N65").Value = 92.878999999999991
$ws.Range("O65").Value = 463.03699999999992
$ws.Range("P65").Value = 412.21200000000005
$ws.Range("Q65").Value = 405.20699999999999
$ws.Range("R65").Value = 416.78100000000001
$ws.Range("S65").Value = 202.08199999999999
$ws.Range("T65").Value = 173.76100000000002
$ws.Range("U65").Value = 29.18
$ws.Range("V65").Value = 43.742000000000004
$ws.Range("W65").Value = 17.265999999999998
$ws.Range("X65").Value = 43.777000000000001
$ws.Range("Y65").Value = 31.625
$ws.Range("Z65").Value = 0

# Update selection to match final cursor position
$ws.Range("F68").Select()
